# Updated cryptos list on Wed Nov  1 14:07:48 UTC 2023 with GitHub Actions
#
# Applies the latest price/volume snapshot to the cryptos worksheet, and
# re-ranks a few coins whose relative order changed (rows 39-42).
#
# NOTE: the "Price" column (D) stores numbers as plain text (e.g.
# "34.610.88", "0.602") so that values like thousand-separated prices
# round-trip exactly. Excel's COM layer auto-detects plain numeric-looking
# strings (e.g. "225.01") and silently converts them to the Number type,
# which would corrupt the cell's stored representation. To keep those
# cells as genuine text (matching the original file) we briefly force the
# cell to Text format before assigning, then restore the default "Normal"
# style afterwards so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Val
    )
    $cell = $ws.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "34.668.74"
$ws.Range("E2").Value = "  +1.09%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "1.807.14"
$ws.Range("E3").Value = "  +0.54%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.11%  "

# --- Row 5: BNB ---
Set-TextValue "D5" "225.01"
$ws.Range("E5").Value = "  -0.74%  "

# --- Row 6: XRP ---
Set-TextValue "D6" "0.603"
$ws.Range("E6").Value = "  +1.03%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  -0.10%  "

# --- Row 8: Solana ---
Set-TextValue "D8" "39.82"
$ws.Range("E8").Value = "  +10.19%  "

# --- Row 9: Cardano ---
Set-TextValue "D9" "0.291"
$ws.Range("E9").Value = "  -0.90%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  -2.18%  "

# --- Row 11: TRON ---
$ws.Range("E11").Value = "  +3.96%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
$ws.Range("D12").Value = "2.067.18"
$ws.Range("E12").Value = "  +0.42%  "

# --- Row 13: WrappedEther ---
$ws.Range("D13").Value = "1.808.32"
$ws.Range("E13").Value = "  +0.44%  "

# --- Row 14: Chainlink ---
Set-TextValue "D14" "10.94"
$ws.Range("E14").Value = "  -2.54%  "

# --- Row 15: Polygon ---
Set-TextValue "D15" "0.636"
$ws.Range("E15").Value = "  -0.48%  "

# --- Row 16: WrappedBTC ---
$ws.Range("D16").Value = "34.683.15"
$ws.Range("E16").Value = "  +1.03%  "

# --- Row 17: Polkadot ---
Set-TextValue "D17" "4.38"
$ws.Range("E17").Value = "  -0.84%  "

# --- Row 18: Litecoin ---
Set-TextValue "D18" "67.90"
$ws.Range("E18").Value = "  -2.27%  "

# --- Row 19: BitcoinCash ---
Set-TextValue "D19" "240.92"
$ws.Range("E19").Value = "  -0.85%  "

# --- Row 20: ShibaInu ---
$ws.Range("E20").Value = "  -1.62%  "

# --- Row 21: Avalanche ---
Set-TextValue "D21" "11.09"
$ws.Range("E21").Value = "  -2.64%  "

# --- Row 22: Dai ---
$ws.Range("E22").Value = "  -0.01%  "

# --- Row 23: Uniswap ---
$ws.Range("E23").Value = "  -1.08%  "

# --- Row 24: Toncoin ---
Set-TextValue "D24" "2.18"
$ws.Range("E24").Value = "  -1.94%  "

# --- Row 25: Monero ---
Set-TextValue "D25" "171.71"
$ws.Range("E25").Value = "  +1.01%  "

# --- Row 26: Cosmos ---
Set-TextValue "D26" "7.69"
$ws.Range("E26").Value = "  -4.69%  "

# --- Row 27: EthereumClassic ---
Set-TextValue "D27" "17.46"
$ws.Range("E27").Value = "  +1.43%  "

# --- Row 28: Stellar ---
$ws.Range("E28").Value = "  +0.26%  "

# --- Row 29: BinanceUSD ---
$ws.Range("E29").Value = "  -0.04%  "

# --- Row 30: PancakeSwap ---
Set-TextValue "D30" "1.23"
$ws.Range("E30").Value = "  -0.94%  "

# --- Row 31: Filecoin ---
Set-TextValue "D31" "3.77"
$ws.Range("E31").Value = "  -0.83%  "

# --- Row 32: Hedera ---
Set-TextValue "D32" "0.0515"
$ws.Range("E32").Value = "  -0.51%  "

# --- Row 33: InternetComputer(DFINITY) ---
Set-TextValue "D33" "3.84"
$ws.Range("E33").Value = "  -2.45%  "

# --- Row 34: LidoDAOToken ---
$ws.Range("E34").Value = "  +1.29%  "

# --- Row 35: ImmutableX ---
$ws.Range("E35").Value = "  -2.05%  "

# --- Row 36: TrustWalletToken ---
$ws.Range("E36").Value = "  +0.30%  "

# --- Row 37: Maker ---
$ws.Range("D37").Value = "1.306.79"
$ws.Range("E37").Value = "  -5.14%  "

# --- Row 38: RenderToken ---
$ws.Range("E38").Value = "  +1.20%  "

# --- Rows 39-42: re-ranked — InjectiveProtocol and Aave moved up past
#     VeChain and WEMIXToken respectively ---

# Row 39: was VeChain -> now InjectiveProtocol
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D39" "14.82"
$ws.Range("E39").Value = "  +12.06%  "

# Row 40: was InjectiveProtocol -> now VeChain
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.0188"
$ws.Range("E40").Value = "  +0.95%  "

# Row 41: was WEMIXToken -> now Aave
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D41" "84.06"
$ws.Range("E41").Value = "  +3.23%  "

# Row 42: was Aave -> now WEMIXToken
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D42" "1.24"
$ws.Range("E42").Value = "  +5.15%  "

# --- Row 43: HuobiToken ---
$ws.Range("E43").Value = "  +0.76%  "

# --- Row 44: MXToken ---
$ws.Range("E44").Value = "  +0.40%  "

# --- Row 45: ARBITRUM ---
Set-TextValue "D45" "0.944"
$ws.Range("E45").Value = "  +0.13%  "

# --- Row 46: Kaspa ---
Set-TextValue "D46" "0.0521"
$ws.Range("E46").Value = "  +4.76%  "

# --- Row 47: RocketPoolETH ---
$ws.Range("D47").Value = "1.966.64"
$ws.Range("E47").Value = "  +0.31%  "

# --- Row 48: FraxShare ---
$ws.Range("E48").Value = "  -2.28%  "

# --- Row 49: PaxDollar ---
$ws.Range("E49").Value = "  -0.05%  "

# --- Row 50: Quant ---
Set-TextValue "D50" "101.47"
$ws.Range("E50").Value = "  -0.90%  "

# --- Row 51: Cronos ---
Set-TextValue "D51" "0.0609"
$ws.Range("E51").Value = "  +0.40%  "
